$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the "Nro Evaluador"/"Nro Direccion" columns out: column G takes what
# used to be in column H, and column H takes what used to be in column J.
# Columns I and J are left blank (formatting kept, contents cleared).
$null = $ws.Range("H1:H8").Copy($ws.Range("G1:G8"))
$null = $ws.Range("J1:J8").Copy($ws.Range("H1:H8"))
$null = $ws.Range("I1:I8").ClearContents()
$null = $ws.Range("J1:J8").ClearContents()

# Fix the RTC certification code typo.
$ws.Range("A6").Value = "RTC1"

# Restore the view: scroll back so column A is visible and leave the
# selection on C9.
$null = $ws.Range("A1").Select()
$null = $ws.Range("C9").Select()
